$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# The START JOB command now documents the optional ARGUMENTS clause.
$ws.Range("A25").Value = "START JOB jobName [ARGUMENTS argName1 argValue1 [, " + [char]0x2026 + "]]]"

# Column A needs to widen (best-fit) to accommodate the longer command text.
$ws.Columns.Item(1).ColumnWidth = 56.1666666666667

# The sheet's AutoFilter was removed.
$ws.AutoFilterMode = $false
